# Insert a new row 2 with "University of Washington" data, pushing the
# existing ranking table (rows 2-25) down to rows 3-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a fresh blank row above the current row 2. This shifts all
#    existing row content (values, formulas, shared-formula groups, and
#    per-cell styles) down by one row.
$ws.Rows("2:2").Insert()

# 2) Populate the new row 2 with the University of Washington entry.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "University of Washington"
$ws.Range("C2").Value = "United States"
$ws.Range("D2").Value = 100
$ws.Range("E2").Formula = "=INT(RAND()*100000)"
$ws.Range("F2").Value = "https://www.washington.edu/"
$ws.Range("G2").Value = "University of Washington.png"

# 3) The Rank (A) and Global Score (D) columns are simple sequential fills
#    that are independent of which university occupies the row, so make
#    sure rows 3-26 keep showing the plain 2..25 / 98..52 sequence (the
#    row-insert above duplicated row 2's old literal values instead of
#    continuing the sequence).
for ($r = 3; $r -le 26; $r++) {
    $ws.Range("A$r").Value = $r - 1
    $ws.Range("D$r").Value = 102 - (2 * ($r - 1))
}

# 4) Hyperlinks are anchored to absolute cell refs and are not moved by
#    Insert(), so rebuild them one row lower than where they used to be.
#    Recreating a hyperlink on a cell re-applies Excel's built-in
#    "Hyperlink" cell style, so restore each cell's original formatting
#    right after.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.topuniversities.com/universities/university-sheffield", "", "", "https://www.topuniversities.com/universities/university-sheffield") | Out-Null
$ws.Range("B3").Style = "Normal"
$ws.Range("B3").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.usnews.com/education/best-global-universities/harvard-university-166027", "", "", "https://www.usnews.com/education/best-global-universities/harvard-university-166027") | Out-Null
$ws.Range("B10").Style = "Normal"
$ws.Range("B10").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.usnews.com/education/best-global-universities/university-of-california-berkeley-110635", "", "", "https://www.usnews.com/education/best-global-universities/university-of-california-berkeley-110635") | Out-Null
$ws.Range("B8").Style = "Normal"
$ws.Range("B8").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.berkeley.edu/") | Out-Null
$ws.Range("F8").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.osu.edu/") | Out-Null
$ws.Range("F6").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F24"), "https://www.mit.edu/") | Out-Null
$ws.Range("F24").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F26"), "https://www.topuniversities.com/universities/university-cambridge") | Out-Null
$ws.Range("F26").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F25"), "https://www.topuniversities.com/universities/university-oxford") | Out-Null
$ws.Range("F25").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.topuniversities.com/universities/stanford-university") | Out-Null
$ws.Range("F13").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.topuniversities.com/universities/eth-zurich-swiss-federal-institute-technology") | Out-Null
$ws.Range("F14").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.topuniversities.com/universities/ucl") | Out-Null
$ws.Range("F15").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.topuniversities.com/universities/university-chicago") | Out-Null
$ws.Range("F17").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.topuniversities.com/universities/university-pennsylvania") | Out-Null
$ws.Range("F18").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.topuniversities.com/universities/cornell-university") | Out-Null
$ws.Range("F7").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.topuniversities.com/universities/university-melbourne") | Out-Null
$ws.Range("F19").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.topuniversities.com/universities/california-institute-technology-caltech") | Out-Null
$ws.Range("F20").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("F23"), "https://www.sydney.edu.au/") | Out-Null
$ws.Range("F23").Style = "Normal"
